$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column F (match formatting of the adjacent header cell)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# time_taken values for rows 2-52 (microsecond-precision timestamps)
$timeTaken = @(
    "2021-10-05 13:41:16.839000",
    "2021-10-05 13:41:16.839012",
    "2021-10-05 13:41:16.839016",
    "2021-10-05 13:41:16.839019",
    "2021-10-05 13:41:16.839022",
    "2021-10-05 13:41:16.839025",
    "2021-10-05 13:41:16.839028",
    "2021-10-05 13:41:16.839031",
    "2021-10-05 13:41:16.839035",
    "2021-10-05 13:41:16.839038",
    "2021-10-05 13:41:16.839041",
    "2021-10-05 13:41:16.839044",
    "2021-10-05 13:41:16.839047",
    "2021-10-05 13:41:16.839050",
    "2021-10-05 13:41:16.839053",
    "2021-10-05 13:41:16.839056",
    "2021-10-05 13:41:16.839059",
    "2021-10-05 13:41:16.839062",
    "2021-10-05 13:41:16.839065",
    "2021-10-05 13:41:16.839068",
    "2021-10-05 13:41:16.839071",
    "2021-10-05 13:41:16.839074",
    "2021-10-05 13:41:16.839077",
    "2021-10-05 13:41:16.839080",
    "2021-10-05 13:41:16.839083",
    "2021-10-05 13:41:16.839086",
    "2021-10-05 13:41:16.839089",
    "2021-10-05 13:41:16.839092",
    "2021-10-05 13:41:16.839095",
    "2021-10-05 13:41:16.839098",
    "2021-10-05 13:41:16.839101",
    "2021-10-05 13:41:16.839104",
    "2021-10-05 13:41:16.839108",
    "2021-10-05 13:41:16.839111",
    "2021-10-05 13:41:16.839114",
    "2021-10-05 13:41:16.839117",
    "2021-10-05 13:41:16.839120",
    "2021-10-05 13:41:16.839123",
    "2021-10-05 13:41:16.839126",
    "2021-10-05 13:41:16.839129",
    "2021-10-05 13:41:16.839132",
    "2021-10-05 13:41:16.839135",
    "2021-10-05 13:41:16.839138",
    "2021-10-05 13:41:16.839142",
    "2021-10-05 13:41:16.839145",
    "2021-10-05 13:41:16.839148",
    "2021-10-05 13:41:16.839151",
    "2021-10-05 13:41:16.839154",
    "2021-10-05 13:41:16.839157",
    "2021-10-05 13:41:16.839160",
    "2021-10-05 13:41:16.839163"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
